# Update countries & provincias Spain
# Applies the 4/5 May 2020 data refresh to the Pais worksheet:
#  - Refresh totals for several countries (Estados Unidos, Alemania, Tunez)
#  - Re-rank Colombia above Noruega/Chequia with its updated numbers
#  - Swap the display order of the tied Belice/Santa Lucia and
#    Burundi/San Cristobal y Nieves rows, moving their data accordingly
#  - Update the "Datos actualizados a ..." timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 00:08"

# Estados Unidos (row 4) - refreshed counts
$ws.Range("B4").Value = 1209703
$ws.Range("C4").Value = 21581
$ws.Range("D4").Value = 184358
$ws.Range("E4").Value = 955869
$ws.Range("G4").Value = 878
$ws.Range("H4").Value = 69476

# Alemania (row 9) - refreshed counts
$ws.Range("B9").Value = 165940
$ws.Range("C9").Value = 276
$ws.Range("E9").Value = 26297
$ws.Range("G9").Value = 77
$ws.Range("H9").Value = 6943

# Colombia moves above Noruega/Chequia with its updated numbers (row 46)
$ws.Range("A46").Value = "Colombia"
$ws.Range("B46").Value = 7973
$ws.Range("C46").Value = 305
$ws.Range("D46").Value = 1807
$ws.Range("E46").Value = 5808
$ws.Range("F46").Value = 122
$ws.Range("G46").Value = 18
$ws.Range("H46").Value = 358

# Noruega shifts down to row 47, keeping its prior figures
$ws.Range("A47").Value = "Noruega"
$ws.Range("B47").Value = 7884
$ws.Range("C47").Value = 37
$ws.Range("D47").Value = 32
$ws.Range("E47").Value = 7638
$ws.Range("F47").Value = 27
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 214

# Chequia shifts down to row 48, keeping its prior figures
$ws.Range("A48").Value = "Chequia"
$ws.Range("B48").Value = 7799
$ws.Range("C48").Value = 18
$ws.Range("D48").Value = 3786
$ws.Range("E48").Value = 3762
$ws.Range("F48").Value = 58
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 251

# Tunez (row 93) - refreshed counts
$ws.Range("B93").Value = 1018
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 406
$ws.Range("E93").Value = 569
$ws.Range("F93").Value = 18
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 43

# Belice now listed before Santa Lucia (row 189 gets Belice's data)
$ws.Range("A189").Value = "Belice"
$ws.Range("B189").Value = 18
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 13
$ws.Range("E189").Value = 3
$ws.Range("F189").Value = 1
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 2

# Santa Lucia shifts down to row 190 with its own data
$ws.Range("A190").Value = "Santa Lucia"
$ws.Range("B190").Value = 18
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 15
$ws.Range("E190").Value = 3
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Burundi now listed before San Cristobal y Nieves (row 198 gets Burundi's data)
$ws.Range("A198").Value = "Burundi"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 7
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

# San Cristobal y Nieves shifts down to row 199 with its own data
$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 8
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0
